$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 276, shifting existing rows 276.. down by one.
$ws.Rows.Item(276).Insert()

# Populate the newly inserted row 276 with the new record.
$ws.Cells.Item(276, 1).Value = 9
$ws.Cells.Item(276, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(276, 3).Value = "Metropolitana"
$ws.Cells.Item(276, 4).Value2 = 45120
$ws.Cells.Item(276, 5).Value = 13
$ws.Cells.Item(276, 6).Value = 100112001
$ws.Cells.Item(276, 7).Value = "Berenjena"
$ws.Cells.Item(276, 8).Value = "Sin especificar"
$ws.Cells.Item(276, 9).Value = "Primera"
$ws.Cells.Item(276, 10).Value = 70
$ws.Cells.Item(276, 11).Value = 7000
$ws.Cells.Item(276, 12).Value = 9000
$ws.Cells.Item(276, 13).Value = 8000
$ws.Cells.Item(276, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(276, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(276, 16).Value = 160
$ws.Cells.Item(276, 17).Value = 50
$ws.Cells.Item(276, 18).Value = "Hortaliza"
